# Update the "last updated" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 09:09"

# Update country statistics (columns: B=Casos totales, C=Nuevos casos,
# D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)

# Row 6 - India
$ws.Cells.Item(6, 2).Value = 850827
$ws.Cells.Item(6, 3).Value = 469
$ws.Cells.Item(6, 4).Value = 536314
$ws.Cells.Item(6, 5).Value = 291817
$ws.Cells.Item(6, 7).Value = 9
$ws.Cells.Item(6, 8).Value = 22696

# Row 48 - Afganistan
$ws.Cells.Item(48, 2).Value = 34451
$ws.Cells.Item(48, 3).Value = 85
$ws.Cells.Item(48, 4).Value = 21216
$ws.Cells.Item(48, 5).Value = 12225
$ws.Cells.Item(48, 7).Value = 16
$ws.Cells.Item(48, 8).Value = 1010

# Row 53 - Armenia
$ws.Cells.Item(53, 2).Value = 31969
$ws.Cells.Item(53, 3).Value = 577
$ws.Cells.Item(53, 4).Value = 19633
$ws.Cells.Item(53, 5).Value = 11771
$ws.Cells.Item(53, 7).Value = 6
$ws.Cells.Item(53, 8).Value = 565

# Row 70 - Uzbekistan
$ws.Cells.Item(70, 2).Value = 12706
$ws.Cells.Item(70, 3).Value = 193
$ws.Cells.Item(70, 5).Value = 4925
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 58

# Row 76 - El Salvador
$ws.Cells.Item(76, 4).Value = 5551
$ws.Cells.Item(76, 5).Value = 3586

# Row 99 - Hungria
$ws.Cells.Item(99, 2).Value = 4234
$ws.Cells.Item(99, 3).Value = 5
$ws.Cells.Item(99, 4).Value = 3036
$ws.Cells.Item(99, 5).Value = 603

# Rows 143-145: Georgia's update pushes it above Uruguay and Zimbabue in the
# ranking, so the three countries occupying ranks 147-149 shift:
#   rank 147 -> Georgia (new data)
#   rank 148 -> Uruguay (previous rank-147 data)
#   rank 149 -> Zimbabue (previous rank-148 data)

# Row 143 (rank 147) - now Georgia
$ws.Cells.Item(143, 1).Value = "Georgia"
$ws.Cells.Item(143, 2).Value = 986
$ws.Cells.Item(143, 3).Value = 5
$ws.Cells.Item(143, 4).Value = 857
$ws.Cells.Item(143, 5).Value = 114
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 15

# Row 144 (rank 148) - now Uruguay
$ws.Cells.Item(144, 1).Value = "Uruguay"
$ws.Cells.Item(144, 2).Value = 986
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 896
$ws.Cells.Item(144, 5).Value = 60
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 30

# Row 145 (rank 149) - now Zimbabue
$ws.Cells.Item(145, 1).Value = "Zimbabue"
$ws.Cells.Item(145, 2).Value = 982
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 320
$ws.Cells.Item(145, 5).Value = 644
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 18
